$wb = $excel.ActiveWorkbook

# Rows in each per-language sheet (7,8,9,10,13,14) correspond to file entries
# that just finished generating their handoff report (Priority becomes "ht"
# and the handoff timestamps get refreshed).
$rows = @(7, 8, 9, 10, 13, 14)

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-08-12 06:25:26"
}

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-08-12 06:25:35"
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-12 06:25:35"
}
